# Update the "Pais" (countries) COVID-19 snapshot sheet:
#  - bump the "last updated" timestamp in A1
#  - refresh totals for Estados Unidos (row 4)
#  - re-rank several countries whose case counts moved them up/down the
#    table (Lituania, Uzbekistan, Madagascar, Uganda, Maldivas, ...),
#    which shifts the country name + numbers shown on several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 08:25"

# Each entry: row, country name (col A), then Casos totales/Nuevos casos/
# Casos activos/Recuperados/Casos criticos/Muertes hoy/Muertes (cols B..H)
$rows = @(
    @(4,   "Estados Unidos",       123776, 198, 3231, 118316, 2666, 8, 2229),
    @(66,  "Lituania",             437,  43,  1,  429, 2, 0, 7),
    @(67,  "Libano",               412,  0,   30, 374, 4, 0, 8),
    @(68,  "Hungria",              408,  65,  34, 361, 6, 2, 13),
    @(69,  "Armenia",              407,  0,   30, 374, 6, 2, 3),
    @(70,  "Marruecos",            402,  0,   12, 365, 1, 0, 25),
    @(72,  "Bulgaria",             338,  7,   11, 320, 8, 0, 7),
    @(97,  "Uzbekistan",           133,  29,  5,  126, 8, 0, 2),
    @(98,  "Senegal",              130,  0,   18, 112, 0, 0, 0),
    @(99,  "Brunei",               120,  0,   25, 94,  1, 0, 1),
    @(100, "Cuba",                 119,  0,   4,  112, 2, 0, 3),
    @(101, "Venezuela",            119,  0,   39, 78,  2, 0, 2),
    @(102, "Sri Lanka",            115,  2,   9,  105, 5, 0, 1),
    @(103, "Honduras",             110,  15,  3,  106, 4, 0, 1),
    @(104, "Afganistan",           110,  0,   2,  104, 0, 0, 4),
    @(128, "Madagascar",           39,   13,  0,  39,  0, 0, 0),
    @(129, "Puerto Rico",          39,   0,   1,  36,  0, 0, 2),
    @(130, "Kenia",                38,   0,   1,  36,  2, 0, 1),
    @(131, "Macao",                34,   0,   10, 24,  0, 0, 0),
    @(132, "Guatemala",            34,   0,   10, 23,  1, 0, 1),
    @(133, "Isla de Man",          32,   0,   0,  32,  0, 0, 0),
    @(134, "Guam",                 32,   0,   0,  31,  0, 0, 1),
    @(135, "Jamaica",              32,   2,   2,  29,  0, 0, 1),
    @(137, "Polinesia Francesa",   30,   0,   0,  30,  0, 0, 0),
    @(138, "Zambia",               28,   0,   0,  28,  0, 0, 0),
    @(139, "Guayana Francesa",     28,   0,   6,  22,  0, 0, 0),
    @(140, "Barbados",             26,   0,   0,  26,  0, 0, 0),
    @(142, "El Salvador",          24,   5,   0,  24,  0, 0, 0),
    @(146, "Maldivas",             17,   1,   11, 6,   0, 0, 0),
    @(147, "Etiopia",              16,   0,   1,  15,  0, 0, 0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    for ($c = 2; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}
